$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Limbo
$ws.Range("A2").Value = 1030
$txt_B_2 = @'
Limbo
'@
$ws.Range("B2").Value = $txt_B_2
$txt_C_2 = @'
Martin Stig Andersen, Jeppe Carlsen, Thomas Krog, Mads Wibroe, Arnt Jensen, Dino Christian Patti, Stine Sørensen, Morten Christian Bramsen
'@
$ws.Range("C2").Value = $txt_C_2
$txt_D_2 = @'
2010-07-21
'@
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = $txt_D_2
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = 4
$txt_F_2 = @'
Limbo is a critically acclaimed 2D puzzle-platformer developed by Playdead, released on July 21, 2010. The game was conceptualized and directed by Arnt Jensen, with significant contributions from developers such as Martin Stig Andersen and Jeppe Carlsen. It was published by Microsoft Studios, Playdead, and 鱼俞 across multiple platforms, including PC, PlayStation, Xbox, and Nintendo Switch. Limbo is celebrated for its minimalist artistic style, which employs a monochrome color palette to create a haunting atmosphere that emphasizes the isolation of its nameless protagonist.
The gameplay of Limbo revolves around navigating a series of puzzles and platforming challenges in a hostile environment. The game is notable for its use of one-hit deaths, which are mitigated by a generous checkpoint system that encourages experimentation and learning. The player must guide the protagonist through various treacherous environments, using simple controls to manipulate objects and solve physics-based puzzles. Despite its black-and-white aesthetics, the game intuitively distinguishes interactive elements from those that pose danger, ensuring a seamless gameplay experience.
Limbo's narrative is sparse and open to interpretation, enhancing its enigmatic charm. The game avoids explicit storytelling, instead using its atmospheric design and sound to evoke emotion and tension without relying on traditional horror tropes. This approach, combined with its abrupt ending, has sparked discussions and analysis, contributing to its status as a significant work in the indie gaming scene. The game's sound design, crafted by Andersen, further immerses players in its eerie world, using ambient sounds to heighten suspense and engagement.
Upon release, Limbo received widespread critical acclaim for its artistic design, innovative gameplay, and unique atmosphere. It earned high praise from various gaming outlets and was recognized with numerous awards, solidifying its position as a milestone in the puzzle-platformer genre. The game's impact extends beyond its immediate success, influencing subsequent indie games with its minimalist approach and atmospheric storytelling.
Limbo's cultural impact is significant, as it demonstrated the potential of indie games to deliver profound and innovative experiences. Its success helped pave the way for other indie developers to explore unconventional narratives and artistic styles, contributing to a broader acceptance of video games as an art form.
'@
$ws.Range("F2").Value = $txt_F_2
$txt_G_2 = @'
['Edge Staff. (2010, July 23). Limbo review. Edge. https://www.edge-online.com/features/limbo-review', "Herold, C. (2010, July 21). 'Limbo' a haunting, memorable game. The New York Times. https://www.nytimes.com/2010/07/22/arts/video-games/limbo.html", 'Orland, K. (2010, July 20). Review: Limbo. Ars Technica. https://arstechnica.com/gaming/2010/07/review-limbo', 'Parkin, S. (2010, July 19). Limbo review. Eurogamer. https://www.eurogamer.net/articles/2010-07-19-limbo-review', 'Totilo, S. (2010, July 27). The Making of Limbo. Kotaku. https://kotaku.com/the-making-of-limbo-5597437']
'@
$ws.Range("G2").Value = $txt_G_2
$txt_H_2 = @'
Average Rating: 4.14
Average Playtime: 3 hours
ESRB Rating: Teen
Metacritic Score: 88
Platforms: PC, Android, PS Vita, PlayStation 4, PlayStation 3, Xbox 360, Linux, macOS, iOS, Xbox One, Nintendo Switch
'@
$ws.Range("H2").Value = $txt_H_2
$txt_I_2 = @'
2025-03-27 18:18:32
'@
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = $txt_I_2
$ws.Range("I2").Style = "Normal"

# Row 3: Terraria
$ws.Range("A3").Value = 422
$txt_B_3 = @'
Terraria
'@
$ws.Range("B3").Value = $txt_B_3
$txt_C_3 = @'
Andrew Spinks, Whitney Baird, Victor Moura, Jim Kjexrud, Jamison Hayes, Scott Lloyd Shelly, Chris Bednarz, Yorai Omer
'@
$ws.Range("C3").Value = $txt_C_3
$txt_D_3 = @'
2011-05-16
'@
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = $txt_D_3
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = 4
$txt_F_3 = @'
Terraria is a critically acclaimed 2D action-adventure sandbox game initially released on May 16, 2011. Developed by Re-Logic, a company founded by Andrew Spinks, the game combines elements of classic action games with the freedom of sandbox-style creativity. It is available on multiple platforms, including PC, macOS, Linux, and a range of consoles and mobile devices such as Xbox One, PlayStation 4, and Nintendo Switch. The game has been published by several companies over its lifespan, including 505 Games and Spike Chunsoft.
The gameplay of Terraria revolves around exploration, crafting, building, and combat. Players start by creating a character and entering a randomly generated 2D world filled with diverse biomes, each with unique resources and enemies. The primary goal is to gather materials to craft increasingly powerful weapons and armor, allowing players to tackle more challenging bosses and events. Players can also construct settlements and invite NPCs to inhabit them, which adds depth to the gameplay by providing access to new items and services. The game supports singleplayer and multiplayer modes, with options for cooperative play.
Terraria's reception has been overwhelmingly positive, praised for its depth, replayability, and the vast array of content available. The pixel-art graphics, combined with a rich soundtrack composed by Scott Lloyd Shelly, contribute to its charming and atmospheric experience. The game's success has led to a thriving modding community, which has further expanded its content and replay value. Terraria has sold millions of copies worldwide, making it one of the best-selling indie games of all time.
Culturally, Terraria has become a significant part of the indie gaming landscape, often compared to games like Minecraft for its creative freedom and sandbox elements. It has influenced a new generation of sandbox and survival games, encouraging developers to explore similar mechanics and styles. The game's ongoing updates and community support have kept it relevant for over a decade, maintaining a dedicated fanbase and ensuring its place in video game history.
'@
$ws.Range("F3").Value = $txt_F_3
$txt_G_3 = @'
['IGN. (2011). Terraria Review. Retrieved from https://www.ign.com/articles/2011/05/23/terraria-review', 'Metacritic. (n.d.). Terraria for PC Reviews. Retrieved from https://www.metacritic.com/game/pc/terraria', 'PC Gamer. (2011). Terraria Review. Retrieved from https://www.pcgamer.com/terraria-review', 'Polygon. (2020). Terraria: Journey’s End is a perfect goodbye to a decade-old game. Retrieved from https://www.polygon.com/reviews/2020/5/21/21266515/terraria-journeys-end-review-final-update-pc', 'Re-Logic. (n.d.). Official Terraria Website. Retrieved from http://www.terraria.org/']
'@
$ws.Range("G3").Value = $txt_G_3
$txt_H_3 = @'
Average Rating: 4.06
Average Playtime: 15 hours
ESRB Rating: Teen
Metacritic Score: 81
Platforms: Xbox 360, Wii U, Nintendo 3DS, Xbox One, PlayStation 4, iOS, PC, macOS, Linux, Nintendo Switch, PlayStation 3, PS Vita, Android
'@
$ws.Range("H3").Value = $txt_H_3
$txt_I_3 = @'
2025-03-27 18:19:04
'@
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = $txt_I_3
$ws.Range("I3").Style = "Normal"

# Row 4: Rocket League
$ws.Range("A4").Value = 3272
$txt_B_4 = @'
Rocket League
'@
$ws.Range("B4").Value = $txt_B_4
$txt_C_4 = @'
Corey Davis, Thomas Silloway, Bobby McCoin, Sarah Hebbler, Jared Cone
'@
$ws.Range("C4").Value = $txt_C_4
$txt_D_4 = @'
2015-07-07
'@
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = $txt_D_4
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = 4
$txt_F_4 = @'
Rocket League is a highly competitive vehicular soccer video game developed and published by Psyonix. Released on July 7, 2015, it is considered a spiritual successor to Psyonix's earlier title, Supersonic Acrobatic Rocket-Powered Battle-Cars. The game quickly gained popularity due to its unique blend of soccer and racing mechanics, offering fast-paced gameplay where players control rocket-powered cars to hit a ball into the opponent's goal. Available on multiple platforms, including Nintendo Switch, Linux, macOS, Xbox One, PC, and PlayStation 4, Rocket League supports various multiplayer modes, including both local and online play, and features cross-platform multiplayer capabilities.
The game's core mechanics revolve around using rocket-boosted cars to compete in soccer matches. Players can choose from different game modes, ranging from casual matches to competitive ranked play, in team formats such as 1v1, 2v2, and 3v3. The game is renowned for its 'easy to learn, hard to master' mechanics, which are complemented by tight controls, allowing for precise maneuvers and strategic gameplay. Boost pads scattered across the field provide the necessary energy to accelerate or launch cars into the air, enabling advanced aerial techniques. Rocket League also offers extensive customization options, allowing players to personalize their vehicles with various decals, paint jobs, and accessories. Collaborations with popular franchises have introduced iconic vehicles like the Batmobile and the DeLorean from Back to the Future.
Upon its release, Rocket League received widespread acclaim from critics and players alike, praised for its innovative gameplay, vibrant graphics, and engaging multiplayer experiences. It was lauded for its ability to appeal to both casual players and competitive gamers, leading to a burgeoning esports scene. The game’s soundtrack, featuring a mix of electronic and indie music, also garnered positive feedback for enhancing the gameplay experience. Rocket League's success is reflected in its high player engagement and the passionate community it has cultivated over the years.
Culturally, Rocket League has had a significant impact on the gaming industry, particularly in the realm of esports. Its inclusion in major esports tournaments and leagues has helped elevate the game's visibility and prestige, solidifying its position as a staple in competitive gaming. The game's cross-platform play feature was a pioneering move that encouraged inclusivity and connectivity among players across different gaming systems. Rocket League's influence extends beyond gaming, as it has become a cultural phenomenon, with its unique concept inspiring a variety of media and merchandise, further embedding it into popular culture.
'@
$ws.Range("F4").Value = $txt_F_4
$txt_G_4 = @'
['IGN Staff. (2015). Rocket League Review. IGN. Retrieved from https://www.ign.com/articles/2015/07/07/rocket-league-review', 'Polygon Staff. (2015). Rocket League is a hit and the new hotness on Twitch. Polygon. Retrieved from https://www.polygon.com/2015/7/10/8929671/rocket-league-sales-twitch-steam', 'GameSpot Staff. (2015). Rocket League Review. GameSpot. Retrieved from https://www.gamespot.com/reviews/rocket-league-review/1900-6416180/', "Tassi, P. (2015). Rocket League's Success Story Is A Tale Of 'Right Place, Right Time'. Forbes. Retrieved from https://www.forbes.com/sites/insertcoin/2015/07/21/rocket-leagues-success-story-is-a-tale-of-right-place-right-time/", 'The Verge Staff. (2016). Rocket League now has more than 19 million players. The Verge. Retrieved from https://www.theverge.com/2016/7/6/12109186/rocket-league-player-count-19-million']
'@
$ws.Range("G4").Value = $txt_G_4
$txt_H_4 = @'
Average Rating: 3.93
Average Playtime: 21 hours
ESRB Rating: Everyone
Metacritic Score: 86
Platforms: Nintendo Switch, Linux, macOS, Xbox One, PC, PlayStation 4
'@
$ws.Range("H4").Value = $txt_H_4
$txt_I_4 = @'
2025-03-27 18:19:20
'@
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = $txt_I_4
$ws.Range("I4").Style = "Normal"

Write-Output "rows added"
